$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the underlying hedge position inputs for row 15 (2 timesteps before exotic maturity)
$ws.Range("L15").Value = -1
$ws.Range("N15").Value = 0.4921
